# Commit: "switch from RunCACounties2 to RunCACounties"
# Applies to the "Data" sheet of the SF.xlsx forecast input workbook:
#   1. Clear the PUI columns (D:F) for rows 5-135 - those figures came from
#      the old RunCACounties2 pull and are no longer populated under the
#      RunCACounties source.
#   2. Populate row 273 (2020-12-16) with the newly-pulled confirmed data.
#   3. Move the frozen-pane (bottomRight) selection down one row, from
#      F273 to F274, reflecting the newly added data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# 1. Clear D5:F135 (values only - keep existing cell styles untouched)
$ws.Range("D5:F135").ClearContents()

# 2. Fill in the newly available row of data
$ws.Range("A273").Value = 44181
$ws.Range("B273").Value = 163
$ws.Range("C273").Value = 9
$ws.Range("D273").Value = 39
$ws.Range("E273").Value = 3
$ws.Range("F273").Value = 172

# 3. Update the view: move the active selection in the frozen bottom-right
#    pane from F273 to F274.
$ws.Activate()
$ws.Range("F274").Select()
